# Actualización definición de hecho
# Se ha establecido una descripción definitiva de hecho para los sprints

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "Código:" paragraph — redefine the testing requirements text.
# ---------------------------------------------------------------------

# a) "...completada, ha sido testeada con tests formales y ha sido revisada."
#    -> "...completada y ha sido revisada. En cuanto a la realización de tests formales,"
$r = $d.Content
$ok = $r.Find.Execute(", ha sido testeada con ", $true, $false, $false, $false, $false, $true, 1, $false, `
    " y ha sido revisada. En cuanto a la realización de ", 2)

# b) "formales y ha sido revisada. Aquel código que no sea una funcionalidad concreta
#     compleja no será necesario que sea testeado formalmente, pero sí serán necesarios"
#    -> "formales, no son obligatorios a no ser que se indiquen explícitamente junto a la
#        tarea. Por otro lado, sí serán obligatorios"
$r = $d.Content
$ok = $r.Find.Execute(" formales y ha sido revisada. Aquel código que no sea una funcionalidad concreta compleja no será necesario que sea testeado formalmente, pero sí serán necesarios ", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    " formales, no son obligatorios a no ser que se indiquen explícitamente junto a la tarea. Por otro lado, sí serán obligatorios ", 2)

# c) "informales (por ejemplo, cuando se producen cambios en el css)." -> "informales."
$r = $d.Content
$ok = $r.Find.Execute(" informales (por ejemplo, cuando se producen cambios en el css).", `
    $true, $false, $false, $false, $false, $true, 1, $false, " informales.", 2)

# d) Split the paragraph right after "...tests informales." — the rest of the
#    paragraph ("No se considerará acabada la tarea...") becomes its own paragraph.
$r = $d.Content
$found = $r.Find.Execute("Por otro lado, sí serán obligatorios tests informales.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $r.Collapse(0)
    $r.InsertParagraphAfter()
}

# ---------------------------------------------------------------------
# 2) Fix missing space / stray run-join: "...lo sepan.En la descripción..."
#    -> "...lo sepan. En la descripción..."
# ---------------------------------------------------------------------
$r = $d.Content
$ok = $r.Find.Execute("lo sepan.En la descripción", $true, $false, $false, $false, $false, $true, 1, $false, `
    "lo sepan. En la descripción", 2)

# ---------------------------------------------------------------------
# 3) Trim six trailing blank paragraphs that sit right before the final
#    tab-stop paragraph near the end of the document (locate it by its
#    own content so earlier paragraph insertions above don't throw the
#    index off).
# ---------------------------------------------------------------------
$paras = $d.Paragraphs
$tabIdx = -1
for ($i = $paras.Count; $i -ge 1; $i--) {
    $p = $paras.Item($i)
    if ($p.Range.Text -eq "`t`r") {
        $tabIdx = $i
        break
    }
}

for ($k = 0; $k -lt 6; $k++) {
    $idx = $tabIdx - 1
    $p = $d.Paragraphs.Item($idx)
    $p.Range.Delete()
    $tabIdx = $tabIdx - 1
}

# ---------------------------------------------------------------------
# 4) Reset the left indent on the final (tab-stop) paragraph to 0.
# ---------------------------------------------------------------------
$tabPara = $d.Paragraphs.Item($tabIdx)
$tabPara.Format.LeftIndent = 0
